$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '33.822.18'
$ws.Range("E2").Value = '  -2.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.776.78'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.67'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '30.98'
$ws.Range("E8").Value = '  -4.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.284'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("E10").Value = '  +5.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.031.94'
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.773.31'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("E14").Value = '  -4.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.624'
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.815.20'
$ws.Range("E16").Value = '  -2.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.21'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.81'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.26'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0773'
$ws.Range("E20").Value = '  +1.29%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.60'
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").Value = '  -3.53%  '
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.38'
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.36'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.49'
$ws.Range("E33").Value = '  -2.28%  '
$ws.Range("E34").Value = '  -2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.392.90'
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.638'
$ws.Range("E36").Value = '  +1.73%  '
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("E39").Value = '  +3.46%  '
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.80'
$ws.Range("E41").Value = '  -4.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  -4.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.11'
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E45").Value = '  -2.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.03'
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.929.24'
$ws.Range("E47").Value = '  -0.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.57'
$ws.Range("E48").Value = '  +1.58%  '
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0119'
$ws.Range("E51").Value = '  -2.51%  '
